# The deck's visible theme ("ppt/theme/theme1.xml", linked from the
# slide master) currently carries the "Integral" colour palette. The
# target edit swaps it for the stock Office "Office Theme" colour
# palette (the font scheme and format/effect scheme are already
# identical between the two themes present in the file, so only the
# twelve theme colours - and, informationally, the scheme name - change).
#
# PowerPoint exposes the live colour scheme through
# Master.ColorScheme.Colors(index).RGB, indices 1-12 mapping to
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (the same order used
# in the OOXML <a:clrScheme> element). RGB values use the classic OLE
# "0x00BBGGRR" long (R + G*256 + B*65536), matching VBA's RGB().

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$cs = $master.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
